# Applies the diff: mark row44 column J with "?", append four new rows
# (45: Intersection of Two Arrays, 46: Intersection of Two Arrays II,
#  47: Longest Mountain in Array, 48: Trapping Rain Water / problem 42)
# and update the sheet view selection to I48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44: add the missing "?" marker in column J -----------------------
$ws.Cells.Item(44, 10).Value = "?"

# Re-use a typed [DateTime] value (rather than the raw serial number) so the
# engine recognises the date and keeps reusing the workbook's existing
# date-formatted style instead of minting a new cellXfs record.
$lastUpdate = [DateTime]"2025-07-03"

# --- Row 45: #349 Intersection of Two Arrays -------------------------------
$ws.Cells.Item(45, 1).Value = 349
$ws.Cells.Item(45, 2).Value = "Intersection of Two Arrays"
$ws.Cells.Item(45, 3).Value = "#hash-table #array #two-pointers #核心 "
$ws.Cells.Item(45, 4).Value = "easy"
$ws.Cells.Item(45, 5).Value = 4
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 5
$ws.Cells.Item(45, 8).Value = $lastUpdate
$ws.Cells.Item(45, 9).Value = $lastUpdate

# --- Row 46: #350 Intersection of Two Arrays II ----------------------------
$ws.Cells.Item(46, 1).Value = 350
$ws.Cells.Item(46, 2).Value = "Intersection of Two Arrays II"
$ws.Cells.Item(46, 3).Value = "#hash-table #array #two-pointers #核心 "
$ws.Cells.Item(46, 4).Value = "easy"
$ws.Cells.Item(46, 5).Value = 4
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 5
$ws.Cells.Item(46, 8).Value = $lastUpdate
$ws.Cells.Item(46, 9).Value = $lastUpdate

# --- Row 47: #845 Longest Mountain in Array --------------------------------
$ws.Cells.Item(47, 1).Value = 845
$ws.Cells.Item(47, 2).Value = "Longest Mountain in Array"
$ws.Cells.Item(47, 3).Value = "#array #two-pointers #核心 "
$ws.Cells.Item(47, 4).Value = "medium"
$ws.Cells.Item(47, 5).Value = 2
$ws.Cells.Item(47, 6).Value = 2
$ws.Cells.Item(47, 7).Value = 20
$ws.Cells.Item(47, 8).Value = $lastUpdate
$ws.Cells.Item(47, 9).Value = $lastUpdate

# --- Row 48: #42 Trapping Rain Water ---------------------------------------
$ws.Cells.Item(48, 1).Value = 42
$ws.Cells.Item(48, 2).Value = "Trapping Rain Water"
$ws.Cells.Item(48, 3).Value = "#array #two-pointers #核心 "
$ws.Cells.Item(48, 4).Value = "hard"
$ws.Cells.Item(48, 5).Value = 1
$ws.Cells.Item(48, 6).Value = 3
$ws.Cells.Item(48, 7).Value = 20
$ws.Cells.Item(48, 8).Value = $lastUpdate
$ws.Cells.Item(48, 9).Value = $lastUpdate
$ws.Cells.Item(48, 10).Value = "?"

# Row heights matching the heights used for similar wrapped-text rows.
$ws.Rows.Item(45).RowHeight = 51
$ws.Rows.Item(46).RowHeight = 51
$ws.Rows.Item(47).RowHeight = 34
$ws.Rows.Item(48).RowHeight = 34

# Update the selection to reflect the new last cell (I48), mirroring the
# workbook's recorded view state after the edit.
$ws.Range("I48").Select() | Out-Null
